$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "2024-11-08 04:30Vykdomi"
$ws.Range("C2").Value = 136.6
$ws.Range("D2").Value = 6.3

$ws.Range("B3").Value = "2024-11-08 04:30Vykdomi"
$ws.Range("C3").Value = 53.2
$ws.Range("D3").Value = 6.8

$ws.Range("B4").Value = "2024-11-08 04:30Vykdomi"
$ws.Range("C4").Value = 44.8
$ws.Range("D4").Value = 7.3

$ws.Range("B5").Value = "2024-11-08 04:30Vykdomi"
$ws.Range("C5").Value = 63
$ws.Range("D5").Value = 7.8

$ws.Range("B6").Value = "2024-11-08 04:30Vykdomi"
$ws.Range("C6").Value = 88.7
$ws.Range("D6").Value = 8

$ws.Range("B7").Value = "2024-11-08 04:30Vykdomi"
$ws.Range("D7").Value = 7.1

$ws.Range("B8").Value = "2024-11-08 04:30Vykdomi"
$ws.Range("C8").Value = 73.2
$ws.Range("D8").Value = 7

$ws.Range("B9").Value = "2024-11-08 04:30Vykdomi"
$ws.Range("C9").Value = 109
$ws.Range("D9").Value = 8.300000000000001

$ws.Range("B10").Value = "2024-11-08 04:30Vykdomi"
$ws.Range("C10").Value = 90.5
$ws.Range("D10").Value = 7.1

$ws.Range("B11").Value = "2024-11-08 04:30Vykdomi"
$ws.Range("C11").Value = 136.7
$ws.Range("D11").Value = 9.800000000000001

$ws.Range("B12").Value = "2024-11-08 04:30Vykdomi"
$ws.Range("D12").Value = 8.199999999999999

$ws.Range("B13").Value = "2024-11-08 04:30Vykdomi"
$ws.Range("C13").Value = 93.59999999999999
$ws.Range("D13").Value = 8.4

$ws.Range("B14").Value = "2024-11-08 04:30Vykdomi"
$ws.Range("C14").Value = 33.1
$ws.Range("D14").Value = 8.199999999999999

$ws.Range("B15").Value = "2024-11-08 04:30Vykdomi"
$ws.Range("C15").Value = 197.9
$ws.Range("D15").Value = 7.3

$ws.Range("B16").Value = "2024-11-08 04:30Vykdomi"
$ws.Range("C16").Value = 12.7
$ws.Range("D16").Value = 6.5

$ws.Range("B17").Value = "2024-11-08 04:30Vykdomi"
$ws.Range("C17").Value = 16
$ws.Range("D17").Value = 7.7

$ws.Range("B18").Value = "2024-11-08 04:30Vykdomi"
$ws.Range("C18").Value = 64.2
$ws.Range("D18").Value = 8.9

$ws.Range("B19").Value = "2024-11-08 04:30Vykdomi"
$ws.Range("D19").Value = 7.7

$ws.Range("B20").Value = "2024-11-08 04:30Vykdomi"
$ws.Range("C20").Value = 71.2

$ws.Range("B21").Value = "2024-11-08 04:30Vykdomi"
$ws.Range("C21").Value = 92.5
$ws.Range("D21").Value = 7.6

$ws.Range("B22").Value = "2024-11-08 04:30Vykdomi"
$ws.Range("C22").Value = 88.8

$ws.Range("B23").Value = "2024-11-08 04:30Vykdomi"
$ws.Range("C23").Value = 138.9

$ws.Range("B24").Value = "2024-11-08 04:30Vykdomi"
$ws.Range("C24").Value = 63.6

$ws.Range("B25").Value = "2024-11-08 04:30Vykdomi"
$ws.Range("C25").Value = 114.8
$ws.Range("D25").Value = 7.6

$ws.Range("B26").Value = "2024-11-08 04:30Vykdomi"
$ws.Range("C26").Value = 125.6
$ws.Range("D26").Value = 6.3

$ws.Range("B27").Value = "2024-11-08 04:30Vykdomi"
$ws.Range("C27").Value = 98.09999999999999
$ws.Range("D27").Value = 5.8

$ws.Range("B28").Value = "2024-11-08 04:30Vykdomi"
$ws.Range("C28").Value = 18
$ws.Range("D28").Value = 6.3

$ws.Range("B29").Value = "2024-11-08 04:30Vykdomi"
$ws.Range("C29").Value = 42.1
$ws.Range("D29").Value = 6.2

$ws.Range("B30").Value = "2024-11-08 04:30Vykdomi"
$ws.Range("C30").Value = 87.8
$ws.Range("D30").Value = 6.7

$ws.Range("B31").Value = "2024-11-08 04:30Vykdomi"
$ws.Range("C31").Value = 300.4
$ws.Range("D31").Value = 7.2

$ws.Range("B32").Value = "2024-11-08 04:30Vykdomi"
$ws.Range("C32").Value = 227.6
$ws.Range("D32").Value = 6.8

$ws.Range("B33").Value = "2024-11-08 04:30Vykdomi"
$ws.Range("C33").Value = 246.8
$ws.Range("D33").Value = 10.5

$ws.Range("B34").Value = "2024-11-08 04:30Vykdomi"
$ws.Range("C34").Value = 31.5

$ws.Range("B35").Value = "2024-11-08 04:30Vykdomi"
$ws.Range("C35").Value = 78.40000000000001
$ws.Range("D35").Value = 7.1

$ws.Range("B36").Value = "2024-11-08 04:30Vykdomi"
$ws.Range("C36").Value = 185.3
$ws.Range("D36").Value = 6.6

$ws.Range("B37").Value = "2024-11-08 04:30Vykdomi"
$ws.Range("C37").Value = 2.3
$ws.Range("D37").Value = 7.2

$ws.Range("B38").Value = "2024-11-08 04:30Vykdomi"
$ws.Range("C38").Value = 20.5
$ws.Range("D38").Value = 7.2

$ws.Range("B39").Value = "2024-11-08 04:30Vykdomi"
$ws.Range("C39").Value = 58
$ws.Range("D39").Value = 6.4

$ws.Range("B40").Value = "2024-11-08 04:30Vykdomi"
$ws.Range("C40").Value = 59.7
$ws.Range("D40").Value = 9

$ws.Range("B41").Value = "2024-11-08 04:30Vykdomi"
$ws.Range("C41").Value = 16.9
$ws.Range("D41").Value = 9

$ws.Range("B42").Value = "2024-11-08 04:30Vykdomi"
$ws.Range("D42").Value = 7.1

$ws.Range("B43").Value = "2024-11-08 04:30Vykdomi"
$ws.Range("C43").Value = 118.4
$ws.Range("D43").Value = 6.9

$ws.Range("B44").Value = "2024-11-08 04:30Vykdomi"
$ws.Range("D44").Value = 7.5

$ws.Range("B45").Value = "2024-11-08 04:30Vykdomi"
$ws.Range("C45").Value = 99.3
$ws.Range("D45").Value = 7

$ws.Range("B46").Value = "2024-11-08 04:30Vykdomi"
$ws.Range("C46").Value = 2.3
$ws.Range("D46").Value = 9.9

$ws.Range("B47").Value = "2024-11-08 04:30Vykdomi"
$ws.Range("C47").Value = 19.6
$ws.Range("D47").Value = 7.9

$ws.Range("B48").Value = "2024-11-08 04:30Vykdomi"
$ws.Range("C48").Value = 19.5
$ws.Range("D48").Value = 6.7

$ws.Range("B49").Value = "2024-11-08 04:30Vykdomi"
$ws.Range("C49").Value = 123.4
$ws.Range("D49").Value = 8.1

$ws.Range("B50").Value = "2024-11-08 04:30Vykdomi"
$ws.Range("C50").Value = 19
$ws.Range("D50").Value = 6.8

$ws.Range("B51").Value = "2024-11-08 04:30Vykdomi"
$ws.Range("C51").Value = 17.3
$ws.Range("D51").Value = 6.7

$ws.Range("B52").Value = "2024-11-08 04:30Vykdomi"
$ws.Range("C52").Value = 95.2
$ws.Range("D52").Value = 7.5

$ws.Range("B53").Value = "2024-11-08 04:30Vykdomi"
$ws.Range("C53").Value = 523.3
$ws.Range("D53").Value = 7.3

$ws.Range("B54").Value = "2024-11-08 04:30Vykdomi"
$ws.Range("C54").Value = 68.2
$ws.Range("D54").Value = 7.8

$ws.Range("B55").Value = "2024-11-08 04:30Vykdomi"
$ws.Range("C55").Value = 141.4
$ws.Range("D55").Value = 7.3

$ws.Range("B56").Value = "2024-11-08 04:30Vykdomi"
$ws.Range("C56").Value = 517.8
$ws.Range("D56").Value = 6.6

$ws.Range("B57").Value = "2024-11-08 04:30Vykdomi"
$ws.Range("C57").Value = 511.7
$ws.Range("D57").Value = 5.2

$ws.Range("B58").Value = "2024-11-08 04:30Vykdomi"
$ws.Range("C58").Value = 101.8
$ws.Range("D58").Value = 8

$ws.Range("B59").Value = "2024-11-08 04:30Vykdomi"
$ws.Range("C59").Value = 140.8
$ws.Range("D59").Value = 9

$ws.Range("B60").Value = "2024-11-08 04:30Vykdomi"
$ws.Range("C60").Value = 95.8
$ws.Range("D60").Value = 8

$ws.Range("B61").Value = "2024-11-08 04:30Vykdomi"
$ws.Range("C61").Value = 136.7
$ws.Range("D61").Value = 7.1

$ws.Range("B62").Value = "2024-11-08 04:30Vykdomi"
$ws.Range("C62").Value = 23.2
$ws.Range("D62").Value = 7.9

$ws.Range("B63").Value = "2024-11-08 04:30Vykdomi"
$ws.Range("C63").Value = 78.40000000000001
$ws.Range("D63").Value = 8.300000000000001
